# Unify the conception of DataNode, DataTable, Entity.
# Rename the generic "PropertyN" worksheet tabs to the standardized
# "DataNode_N" naming convention, and leave the second sheet
# (DataNode_2) as the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "DataNode_1"
$ws2.Name = "DataNode_2"

$ws2.Activate()
